# Insert a new data row at row 653 ("Femacal de La Calera" / Apio price
# series), shifting the existing rows 653:727 down to 654:728, then
# populate the newly-inserted row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("653:653").Insert()

$ws.Range("A653").Value = 3
$ws.Range("B653").Value = 'Femacal de La Calera'
$ws.Range("C653").Value = 'Coquimbo'
$ws.Range("D653").Value = 45212
$ws.Range("E653").Value = 5
$ws.Range("F653").Value = 100112017
$ws.Range("G653").Value = 'Apio'
$ws.Range("H653").Value = 'Americana (o)'
$ws.Range("I653").Value = 'Primera'
$ws.Range("J653").Value = 130
$ws.Range("K653").Value = 7000
$ws.Range("L653").Value = 7000
$ws.Range("M653").Value = 7000
$ws.Range("N653").Value = '$/docena de matas'
$ws.Range("O653").Value = 'Pan de Azúcar'
$ws.Range("P653").Value = 1167
$ws.Range("Q653").Value = 6
$ws.Range("R653").Value = 'Hortaliza'
